$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) RECAP sheet: unit of measure row for Cost changes from "USD/kW" to
#    "MUSD/kW" and a new column C is added with "(million USD)".
# ---------------------------------------------------------------------------
$recap = $wb.Worksheets.Item("RECAP")

# Pick up the look of an existing "plain Arial 10" labelled cell (A28) so the
# new/changed cells match the rest of the sheet's label formatting.
$recap.Range("A28").Copy() | Out-Null
$recap.Range("B37:C37").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$recap.Range("B37").Value = "MUSD/kW"
$recap.Range("C37").Value = "(million USD)"
$recap.Range("B37:C37").Font.Name = "Arial"
$recap.Range("B37:C37").Font.Size = 10
$recap.Range("B37:C37").Font.Bold = $false

# ---------------------------------------------------------------------------
# 2) Cost_Onshore / Cost_Offshore: re-run of results - the average investment
#    cost row (row 2) is now expressed in MUSD/kW instead of USD/kW, so every
#    value is rescaled by 1,000,000.
# ---------------------------------------------------------------------------
$costOnshore = $wb.Worksheets.Item("Cost_Onshore")
for ($col = 1; $col -le 101; $col++) {
    $cell = $costOnshore.Cells.Item(2, $col)
    $cell.Value = $cell.Value() / 1000000
}

$costOffshore = $wb.Worksheets.Item("Cost_Offshore")
for ($col = 1; $col -le 101; $col++) {
    $cell = $costOffshore.Cells.Item(2, $col)
    $cell.Value = $cell.Value() / 1000000
}

# ---------------------------------------------------------------------------
# 3) Restore/update the selection on each touched sheet, finishing on
#    Cost_Offshore so it becomes the active tab (matching the new
#    activeTab in the saved workbook view).
# ---------------------------------------------------------------------------
$recap.Activate()
$recap.Range("C38").Select()

$costOnshore.Activate()
$costOnshore.Range("AY9").Select()

$costOffshore.Activate()
$costOffshore.Range("W16").Select()
